$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update data rows 2-5 with new sensor readings
# Row 2
$ws.Cells.Item(2,1).Value = 45032.50694444445
$ws.Cells.Item(2,2).Value = 4.928
$ws.Cells.Item(2,3).Value = 5.344
$ws.Cells.Item(2,4).Value = 0
$ws.Cells.Item(2,5).Value = 5.978
$ws.Cells.Item(2,6).Value = 10.253
$ws.Cells.Item(2,7).Value = 2.511
$ws.Cells.Item(2,8).Value = 8.567
$ws.Cells.Item(2,9).Value = 4.041
$ws.Cells.Item(2,10).Value = 1.748
$ws.Cells.Item(2,11).Value = 4.53
$ws.Cells.Item(2,12).Value = 5.048
$ws.Cells.Item(2,13).Value = 5.182
$ws.Cells.Item(2,14).Value = 0.791
$ws.Cells.Item(2,15).Value = 3.479
$ws.Cells.Item(2,16).Value = 3.593
$ws.Cells.Item(2,17).Value = 1.503
$ws.Cells.Item(2,18).Value = 1.042
$ws.Cells.Item(2,19).Value = 0.417
$ws.Cells.Item(2,20).Value = 43.098
$ws.Cells.Item(2,21).Value = 7.946
$ws.Cells.Item(2,22).Value = 4.36
$ws.Cells.Item(2,23).Value = 6.645
$ws.Cells.Item(2,24).Value = 2.748
$ws.Cells.Item(2,25).Value = 0.484
$ws.Cells.Item(2,26).Value = 2.765
$ws.Cells.Item(2,27).Value = 1.536
$ws.Cells.Item(2,28).Value = 3.068
$ws.Cells.Item(2,29).Value = 2.75
$ws.Cells.Item(2,30).Value = 5.612
$ws.Cells.Item(2,31).Value = 0
$ws.Cells.Item(2,32).Value = 5.035
$ws.Cells.Item(2,33).Value = 2.082
$ws.Cells.Item(2,34).Value = 3.31

# Row 3
$ws.Cells.Item(3,1).Value = 45032.51388888889
$ws.Cells.Item(3,2).Value = 0.527
$ws.Cells.Item(3,3).Value = 1.252
$ws.Cells.Item(3,4).Value = 0.056
$ws.Cells.Item(3,5).Value = 0.245
$ws.Cells.Item(3,6).Value = 1.253
$ws.Cells.Item(3,7).Value = 0
$ws.Cells.Item(3,8).Value = 3.429
$ws.Cells.Item(3,9).Value = 0.257
$ws.Cells.Item(3,10).Value = 0.145
$ws.Cells.Item(3,11).Value = 0.731
$ws.Cells.Item(3,12).Value = 0.903
$ws.Cells.Item(3,13).Value = 0.149
$ws.Cells.Item(3,14).Value = 0
$ws.Cells.Item(3,15).Value = 0.305
$ws.Cells.Item(3,16).Value = 0
$ws.Cells.Item(3,17).Value = 0.094
$ws.Cells.Item(3,18).Value = 0.458
$ws.Cells.Item(3,19).Value = 0.243
$ws.Cells.Item(3,20).Value = 0
$ws.Cells.Item(3,21).Value = 0.803
$ws.Cells.Item(3,22).Value = 0.844
$ws.Cells.Item(3,23).Value = 1.411
$ws.Cells.Item(3,24).Value = 0.342
$ws.Cells.Item(3,25).Value = 0.063
$ws.Cells.Item(3,26).Value = 1.076
$ws.Cells.Item(3,27).Value = 0.146
$ws.Cells.Item(3,28).Value = 0.4
$ws.Cells.Item(3,29).Value = 0.223
$ws.Cells.Item(3,30).Value = 1.035
$ws.Cells.Item(3,31).Value = 0
$ws.Cells.Item(3,32).Value = 2.674
$ws.Cells.Item(3,33).Value = 0.17
$ws.Cells.Item(3,34).Value = 0.14

# Row 4
$ws.Cells.Item(4,1).Value = 45032.52083333334
$ws.Cells.Item(4,2).Value = 13.05
$ws.Cells.Item(4,3).Value = 10.322
$ws.Cells.Item(4,4).Value = 0.458
$ws.Cells.Item(4,5).Value = 27.605
$ws.Cells.Item(4,6).Value = 23.597
$ws.Cells.Item(4,7).Value = 10.411
$ws.Cells.Item(4,8).Value = 33.233
$ws.Cells.Item(4,9).Value = 15.481
$ws.Cells.Item(4,10).Value = 6.899
$ws.Cells.Item(4,11).Value = 10.822
$ws.Cells.Item(4,12).Value = 11.627
$ws.Cells.Item(4,13).Value = 11.766
$ws.Cells.Item(4,14).Value = 3.173
$ws.Cells.Item(4,15).Value = 10.078
$ws.Cells.Item(4,16).Value = 14.004
$ws.Cells.Item(4,17).Value = 8.214
$ws.Cells.Item(4,18).Value = 0.439
$ws.Cells.Item(4,19).Value = 0.556
$ws.Cells.Item(4,20).Value = 147.599
$ws.Cells.Item(4,21).Value = 27.766
$ws.Cells.Item(4,22).Value = 9.667999999999999
$ws.Cells.Item(4,23).Value = 19.284
$ws.Cells.Item(4,24).Value = 9.988
$ws.Cells.Item(4,25).Value = 1.337
$ws.Cells.Item(4,26).Value = 16.471
$ws.Cells.Item(4,27).Value = 8.147
$ws.Cells.Item(4,28).Value = 7.344
$ws.Cells.Item(4,29).Value = 8.506
$ws.Cells.Item(4,30).Value = 12.313
$ws.Cells.Item(4,31).Value = 0
$ws.Cells.Item(4,32).Value = 29.315
$ws.Cells.Item(4,33).Value = 5.332
$ws.Cells.Item(4,34).Value = 11.501

# Row 5
$ws.Cells.Item(5,1).Value = 45032.52777777778
$ws.Cells.Item(5,2).Value = 19.78
$ws.Cells.Item(5,3).Value = 15.23
$ws.Cells.Item(5,4).Value = 0.68
$ws.Cells.Item(5,5).Value = 42.35
$ws.Cells.Item(5,6).Value = 35.64
$ws.Cells.Item(5,7).Value = 15.73
$ws.Cells.Item(5,8).Value = 58.16
$ws.Cells.Item(5,9).Value = 23.68
$ws.Cells.Item(5,10).Value = 10.64
$ws.Cells.Item(5,11).Value = 16.24
$ws.Cells.Item(5,12).Value = 17.42
$ws.Cells.Item(5,13).Value = 18.03
$ws.Cells.Item(5,14).Value = 4.9
$ws.Cells.Item(5,15).Value = 15.35
$ws.Cells.Item(5,16).Value = 21.7
$ws.Cells.Item(5,17).Value = 12.6
$ws.Cells.Item(5,18).Value = 0.43
$ws.Cells.Item(5,19).Value = 0.72
$ws.Cells.Item(5,20).Value = 227.58
$ws.Cells.Item(5,21).Value = 42.78
$ws.Cells.Item(5,22).Value = 14.44
$ws.Cells.Item(5,23).Value = 29.37
$ws.Cells.Item(5,24).Value = 15.29
$ws.Cells.Item(5,25).Value = 2.02
$ws.Cells.Item(5,26).Value = 28.41
$ws.Cells.Item(5,27).Value = 12.46
$ws.Cells.Item(5,28).Value = 11.09
$ws.Cells.Item(5,29).Value = 12.95
$ws.Cells.Item(5,30).Value = 18.39
$ws.Cells.Item(5,31).Value = 0
$ws.Cells.Item(5,32).Value = 52.55
$ws.Cells.Item(5,33).Value = 8.109999999999999
$ws.Cells.Item(5,34).Value = 17.63

# 2) Remove the old row 6 (data now ends at row 5)
$ws.Rows(6).Delete()

# 3) Adjust column widths to match new layout
$offset = 0.8333333333333334
$ws.Columns(2).ColumnWidth = 7 - $offset
$ws.Columns(3).ColumnWidth = 8 - $offset
$ws.Columns(7).ColumnWidth = 8 - $offset
$ws.Columns(11).ColumnWidth = 8 - $offset
$ws.Columns(15).ColumnWidth = 8 - $offset
$ws.Columns(31).ColumnWidth = 5 - $offset
$ws.Columns(32).ColumnWidth = 8 - $offset
